$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1079986666666667
$ws.Range("H2").Value = 0.323996
$ws.Range("I2").Value = 0.004187739561209694
$ws.Range("J2").Value = 0.004187739561209694
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 6.592412095191557
$ws.Range("R2").Value = 59.33170885672401
$ws.Range("S2").Value = 0.0008558120936668625
$ws.Range("T2").Value = 0.0008558120936668625
$ws.Range("G3").Value = 0.1079986666666667
$ws.Range("H3").Value = 0.323996
$ws.Range("I3").Value = 0.004187739561209694
$ws.Range("J3").Value = 0.004187739561209694
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 11.48182057537867
$ws.Range("R3").Value = 103.336385178408
$ws.Range("S3").Value = 0.001490544092789538
$ws.Range("T3").Value = 0.001490544092789538
$ws.Range("G4").Value = 0.1079986666666667
$ws.Range("H4").Value = 0.323996
$ws.Range("I4").Value = 0.004187739561209694
$ws.Range("J4").Value = 0.004187739561209694
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 14.18437309012089
$ws.Range("R4").Value = 127.659357811088
$ws.Range("S4").Value = 0.001841383374753294
$ws.Range("T4").Value = 0.001841383374753294
$ws.Range("I5").Value = 0.9687110856121154
$ws.Range("J5").Value = 0.9687110856121155
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 1524.961756621445
$ws.Range("R5").Value = 13724.655809593
$ws.Range("S5").Value = 0.1979671013964689
$ws.Range("T5").Value = 0.1979671013964689
$ws.Range("I6").Value = 0.9687110856121154
$ws.Range("J6").Value = 0.9687110856121155
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.3447937879550905
$ws.Range("T6").Value = 0.3447937879550906
$ws.Range("I7").Value = 0.9687110856121154
$ws.Range("J7").Value = 0.9687110856121155
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 3281.139921435111
$ws.Range("R7").Value = 29530.259292916
$ws.Range("S7").Value = 0.425950196260556
$ws.Range("T7").Value = 0.4259501962605561
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6989190000000001
$ws.Range("H8").Value = 2.096757
$ws.Range("I8").Value = 0.02710117482667488
$ws.Range("J8").Value = 0.02710117482667489
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 42.66313845688701
$ws.Range("R8").Value = 383.9682461119831
$ws.Range("S8").Value = 0.005538432567317651
$ws.Range("T8").Value = 0.005538432567317652
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6989190000000001
$ws.Range("H9").Value = 2.096757
$ws.Range("I9").Value = 0.02710117482667488
$ws.Range("J9").Value = 0.02710117482667489
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 74.305200262254
$ws.Range("R9").Value = 668.7468023602861
$ws.Range("S9").Value = 0.009646133780556282
$ws.Range("T9").Value = 0.009646133780556286
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6989190000000001
$ws.Range("H10").Value = 2.096757
$ws.Range("I10").Value = 0.02710117482667488
$ws.Range("J10").Value = 0.02710117482667489
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 91.794909712844
$ws.Range("R10").Value = 826.1541874155961
$ws.Range("S10").Value = 0.01191660847880095
$ws.Range("T10").Value = 0.01191660847880095
